$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Helper: split a (already correctly-texted) cell's single run into
# several runs of the given lengths. All runs keep identical character
# formatting - this mirrors how the target document represents the new
# dates as separate <w:r> elements even though their rPr is identical.
# Splitting is done by toggling a character-formatting property (Bold
# on, then back off) over a growing sub-range; that forces a run
# boundary without changing the final visible formatting.
# -----------------------------------------------------------------------
function Split-CellRuns {
    param($Cell, [int[]]$RunLengths)

    $after = $Cell.Range
    $start = $after.Start
    $end = $after.End

    $running = $start
    for ($i = 0; $i -lt $RunLengths.Length - 1; $i++) {
        $running = $running + $RunLengths[$i]
        if ($running -lt $end) {
            $splitRange = $d.Range($start, $running)
            $splitRange.Font.Bold = 1
            $splitRange2 = $d.Range($start, $running)
            $splitRange2.Font.Bold = 0
        }
    }
}

$table2 = $d.Tables.Item(2)
$table3 = $d.Tables.Item(3)

# --- Table 2 (schedule): Task / Time Duration -------------------------
# These "old" strings (date – date) are unique across the whole document,
# so Find/Replace (ReplaceAll) is safe to use here; it also collapses the
# match to a single run which Split-CellRuns then re-splits.

# Row 2: Test Plan Creation -> "28/05/2023 – 30/05/2023" -> "09/06/2023 – 11/06/2023"
$old1 = "28/05/2023" + [char]32 + [char]8211 + [char]32 + "30/05/2023"
$new1 = "09/06/2023" + [char]32 + [char]8211 + [char]32 + "11/06/2023"
$cellA = $table2.Cell(2, 2)
$cellA.Range.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 0, $false, $new1, 2) | Out-Null
Split-CellRuns $table2.Cell(2, 2) @(2, 2, 1, 5, 3, 2, 2, 1, 5)

# Row 3: Test Scenarios and Test Case Creation -> "31/05/2023 – 03/06/2023" -> "11/06/2023 – 22/06/2023"
$old2 = "31/05/2023" + [char]32 + [char]8211 + [char]32 + "03/06/2023"
$new2 = "11/06/2023" + [char]32 + [char]8211 + [char]32 + "22/06/2023"
$cellB = $table2.Cell(3, 2)
$cellB.Range.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 0, $false, $new2, 2) | Out-Null
Split-CellRuns $table2.Cell(3, 2) @(2, 2, 1, 6, 1, 1, 2, 8)

# Row 4: Test Case Execution / Summary Reports Submission -> "03/06/2023 – 04/06/2023" -> "22/06/2023 – 23/06/2023"
$old3 = "03/06/2023" + [char]32 + [char]8211 + [char]32 + "04/06/2023"
$new3 = "22/06/2023" + [char]32 + [char]8211 + [char]32 + "23/06/2023"
$cellC = $table2.Cell(4, 2)
$cellC.Range.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 0, $false, $new3, 2) | Out-Null
Split-CellRuns $table2.Cell(4, 2) @(2, 11, 2, 8)

# --- Table 3 (deliverables): Deliverable / Description / Target Completion Date --
# These single-date cells are NOT unique strings document-wide, so Find
# cannot be safely used (this runtime's Find always scans from the start
# of the document, regardless of which Range invoked it). Instead, set
# the cell Range's .Text directly - that mutation is properly scoped to
# the exact Range's character offsets.

# Row 2: Test Plan -> "30/05/2023" -> "11/06/2023"
$table3.Cell(2, 3).Range.Text = "11/06/2023"
Split-CellRuns $table3.Cell(2, 3) @(2, 2, 1, 5)

# Row 3: Test Cases -> "03/06/2023" -> "22/06/2023"
$table3.Cell(3, 3).Range.Text = "22/06/2023"
Split-CellRuns $table3.Cell(3, 3) @(2, 8)

# Row 4: Bug Reports -> "04/06/2023" -> "23/06/2023"
$table3.Cell(4, 3).Range.Text = "23/06/2023"
Split-CellRuns $table3.Cell(4, 3) @(2, 8)

# Row 5: Summary Report -> "04/06/2023" -> "23/06/2023"
$table3.Cell(5, 3).Range.Text = "23/06/2023"
Split-CellRuns $table3.Cell(5, 3) @(2, 8)

Write-Host "Done."
